$wb = $excel.ActiveWorkbook

$wsFull = $wb.Worksheets.Item("Full results")
$wsPlot = $wb.Worksheets.Item("For plotting")

# --- Sheet "Full results" ---
$wsFull.Range("H2").Value = 0.590090404453821
$wsFull.Range("I2").Value = 0.295582270908201
$wsFull.Range("O2").Value = 0.410178944660381

$wsFull.Range("F3").Value = 0.581023631753013
$wsFull.Range("G3").Value = 0.315306377745138

$wsFull.Range("C4").Value = 0.60493272865191
$wsFull.Range("D4").Value = 0.395523933953698
$wsFull.Range("E4").Value = 1.00045666260561
$wsFull.Range("J4").Value = 0.395343395207135
$wsFull.Range("K4").Value = 0.315162454701054
$wsFull.Range("L4").Value = -0.0090626340726711
$wsFull.Range("M4").Value = 0.0148355494532462
$wsFull.Range("N4").Value = 0.306099820628382

# --- Sheet "For plotting" ---
$wsPlot.Range("C2").Value = 0.395343395207135
$wsPlot.Range("D2").Value = 0.327346032413468
$wsPlot.Range("E2").Value = 0.463340758000802

$wsPlot.Range("C3").Value = 0.306099820628382
$wsPlot.Range("D3").Value = 0.23252778120793
$wsPlot.Range("E3").Value = 0.379671860048835

$wsPlot.Range("C4").Value = 0.410178944660381
$wsPlot.Range("D4").Value = 0.338858924177425
$wsPlot.Range("E4").Value = 0.481498965143336
